$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Lista de archivos objetados" - 5 new data rows appended to the existing
# CUE/N serv/N-E/HWID/BT/Motivo table (rows 174-178).

$data = @(
    @{ A = 60793500; B = "00000670438a263"; C = "6650A9C5AD9F19627421"; CBold = $false; D = "7427EA97CE37";       E = 0;    F = "Revisar n° de servidor, de ser correcto adjuntar llave pública" },
    @{ A = 60793500; B = "00000670438a263"; C = "EE183CE07CFBD86BF819"; CBold = $true;  D = "C03FD5287B07";        E = 1;    F = "Revisar n° de servidor, de ser correcto adjuntar llave pública" },
    @{ A = 60793500; B = "00000670438a263"; C = "61E48A760302DF984B08"; CBold = $false; D = "C03FD5177126";        E = 4;    F = "Revisar n° de servidor, de ser correcto adjuntar llave pública" },
    @{ A = 60793500; B = "00000670438a263"; C = "EE183CE07CFBD86BF819"; CBold = $true;  D = "7427EA9082A8";        E = "0b"; F = "Revisar n° de servidor, de ser correcto adjuntar llave pública" },
    @{ A = 60793500; B = "00000670438a263"; C = "EE183CE07CFBD86BF819"; CBold = $true;  D = "C03FD52891A8";        E = 1;    F = "Revisar n° de servidor, de ser correcto adjuntar llave pública" }
)

$startRow = 174
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F

    # A:F -> text format, wrapped, general horizontal alignment (matches the
    # rest of the "objetados" rows above them).
    $rngAF = $ws.Range("A" + $r + ":F" + $r)
    $rngAF.NumberFormat = "@"
    $rngAF.WrapText = $true
    $rngAF.HorizontalAlignment = 1

    $cC = $ws.Cells.Item($r, 3)
    $cC.Font.Bold = $row.CBold

    # G/H/I carry the same trailing (empty) formatting as the columns to the
    # right of the table elsewhere in the sheet.
    $ws.Cells.Item($r, 7).NumberFormat = "@"

    $ws.Cells.Item($r, 8).NumberFormat = "General"

    $cI = $ws.Cells.Item($r, 9)
    $cI.NumberFormat = "@"
    $cI.HorizontalAlignment = -4131
}

Write-Host "done"
